$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.518.26"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.914.55"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4812"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06709"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "110.22"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.14"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("D12").Value = "1.912.65"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07561"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.243"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6672"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "289.71"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "30.529.31"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007581"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.88"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "2.167.06"
$ws.Range("E21").Value = "  +1.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.468"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.390"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.434"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.48"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.138"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1064"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.154"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.018"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04984"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7264"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.743"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02049"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.007"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4403"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8634"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.889"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "68.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.293"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.307"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1238"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("E51").Value = "  +2.94%  "
